$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename + update version string ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Self-assessment checklist"
$ws1.Range("A1").Value = "Onderstaande checklist kan gebruikt worden voor het uitvoeren van een assessment tegen de Kwaliteitsaanpak ICTU Software Realisatie versie 2.0-build.0, 14-08-2019."

# --- Sheet 2: duplicate sheet1 (to inherit exact column widths / styles),
#     then strip it down to the 2-row "verbeteracties" (action list) sheet ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Self-assessment verbeteracties"

# Remove comments copied along with the sheet (iterate from the front since
# the collection shrinks as items are deleted)
while ($ws2.Comments.Count -gt 0) {
  $ws2.Comments.Item(1).Delete()
}

# Drop all rows except the two header rows
$ws2.Rows("3:76").Delete()

# Drop conditional formatting / data validation inherited from sheet1
$ws2.Cells.FormatConditions.Delete()
$ws2.Cells.Validation.Delete()

# Remove the frozen panes inherited from sheet1
$ws2.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Activate()

# --- Populate the new header content ---
$ws2.Range("A1").Value = "Onderstaande actielijst kan gebruikt worden om acties n.a.v. de self-assessment bij te houden."
$ws2.Range("A2").Value = "Datum"
$ws2.Range("B2").Value = "Actie"
$ws2.Range("C2").Value = "Status"
$ws2.Range("D2").Value = "Toelichting"

$ws1.Activate()
